$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (shifts rows 13..23 down to 14..24,
# and correctly shifts the existing row heights with them).
$ws.Rows.Item(13).Insert()

# The insert copies column A's bold style down into the new A13 even
# though that row should have no entry in column A at all - strip it
# back to a truly blank, unstyled cell.
$ws.Range("A13").Style = "Normal"
$ws.Range("A13").Value = ""

# B13/C13 are brand new cells; give them the same look (wrapped text)
# as the rest of column B/C before filling them in, since a fresh
# cell otherwise inherits column A's bold style.
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)

# --- Fix values that were shifted by one row in the original sheet ---

# Row 10 (Objetivos:) should hold the Portuguese objectives text, not the
# professor's name.
$ws.Range("B10").Value = "Fornecer aos alunos os conceitos básicos e técnicas de dimensionamento dos principais processos e operações unitárias envolvidas no escoamento de fluidos, sistemas particulados e troca térmica."
$ws.Range("C10").Value = "Fornecer aos alunos os conceitos básicos e técnicas de dimensionamento dos principais processos e operações unitárias envolvidas no escoamento de fluidos, sistemas particulados e troca térmica."

# New row 13 (Docentes responsáveis: label is in A12) gets the professor's name.
$ws.Range("B13").Value = "4780627 - Ana Lucia Gabas Ferreira"
$ws.Range("C13").Value = "4780627 - Ana Lucia Gabas Ferreira"

# Row 14 (Programa resumido:) should hold the short syllabus text in
# Portuguese, not "Semestral".
$ws.Range("B14").Value = "Operações unitárias e processos: reologia de fluidos, dimensionamento de tubulações e acessórios, bombeamento, agitação e mistura, caracterização de partículas e leito de partículas, sedimentação, filtração, processos com membranas. Operações unitárias de troca térmica: trocadores de calor e evaporadores."
$ws.Range("C14").Value = "Operações unitárias e processos: reologia de fluidos, dimensionamento de tubulações e acessórios, bombeamento, agitação e mistura, caracterização de partículas e leito de partículas, sedimentação, filtração, processos com membranas. Operações unitárias de troca térmica: trocadores de calor e evaporadores."

# Row 16 (Programa:) should hold the detailed Portuguese syllabus, not a date.
$ws.Range("B16").Value = "- Reologia de fluidos,- Dimensionamento de tubulações,- Acessórios e bombeamento para fluidos industriais,- Agitação e mistura,- Caracterização de partículas e leito de partículas,- Sedimentação,- Filtração,- Processos com membranas.- Operações unitárias de troca térmica: trocadores de calor e evaporadores."
$ws.Range("C16").Value = "- Reologia de fluidos,- Dimensionamento de tubulações,- Acessórios e bombeamento para fluidos industriais,- Agitação e mistura,- Caracterização de partículas e leito de partículas,- Sedimentação,- Filtração,- Processos com membranas.- Operações unitárias de troca térmica: trocadores de calor e evaporadores."

# Row 19 (Critério:) should hold the evaluation method text.
$ws.Range("B19").Value = "Avaliação composta por duas provas."
$ws.Range("C19").Value = "Avaliação composta por duas provas."

# Row 20 (Norma de recuperação:) should hold the grade-averaging criterion.
$ws.Range("B20").Value = "Média das notas das provas."
$ws.Range("C20").Value = "Média das notas das provas."

# Row 21 (Bibliografia:) should hold the recovery-exam rule text.
$ws.Range("B21").Value = "Prova única com todo o conteúdo da disciplina, sendo que a nota [(nota final do semestre + nota de recuperação)/2] deverá ser igual ou superior a 5,0 (cinco)."
$ws.Range("C21").Value = "Prova única com todo o conteúdo da disciplina, sendo que a nota [(nota final do semestre + nota de recuperação)/2] deverá ser igual ou superior a 5,0 (cinco)."

# Row 22 (Requisitos: label lives in A23 now) should hold the bibliography text.
$ws.Range("B22").Value = "Bibliografia básica:DI BERNARDO, L., Métodos e Técnicas de Tratamento de Água, ABES, Rio de Janeiro, Brasil, 1992.FOUST, A.S., WENZEL, L. A., CLUMP, C.W., MAUS, L., ANDERSEN, L.B. Princípio das operações unitárias. Rio de Janeiro: Editora Guanabara Dois, 1982.GEANKOPLIS, C.J. Procesos de transporte y operaciones unitarias. Compañía Editorial Continental, S.A. de C.V. México, D.F., 1998.PERRY, R.H. and CHILTON, C.H. Manual de Engenharia Química. 5a ed., Guanabara Dois, Rio de Janeiro, 1986.REYNOLDS, T.D.; RICHARDS, P. Unit Operations and Processes in environmental Engineering. PWS Publishing, 1995.MACINTYRE, A.J. Bombas e Instalações de Bombeamento. LTC, Rio de Janeiro, 1997"
$ws.Range("C22").Value = "Bibliografia básica:DI BERNARDO, L., Métodos e Técnicas de Tratamento de Água, ABES, Rio de Janeiro, Brasil, 1992.FOUST, A.S., WENZEL, L. A., CLUMP, C.W., MAUS, L., ANDERSEN, L.B. Princípio das operações unitárias. Rio de Janeiro: Editora Guanabara Dois, 1982.GEANKOPLIS, C.J. Procesos de transporte y operaciones unitarias. Compañía Editorial Continental, S.A. de C.V. México, D.F., 1998.PERRY, R.H. and CHILTON, C.H. Manual de Engenharia Química. 5a ed., Guanabara Dois, Rio de Janeiro, 1986.REYNOLDS, T.D.; RICHARDS, P. Unit Operations and Processes in environmental Engineering. PWS Publishing, 1995.MACINTYRE, A.J. Bombas e Instalações de Bombeamento. LTC, Rio de Janeiro, 1997"
